$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.049.03"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.834.53"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'242.86"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'0.6280"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07589"
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("D9").Value = "'0.2928"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'22.64"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.07743"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "1.846.34"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'4.967"
$ws.Range("D14").Value = "'0.6658"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'0.00001010"
$ws.Range("E15").Value = "  +16.65%  "
$ws.Range("D16").Value = "'83.22"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "'6.076"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "29.062.70"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'227.25"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "'12.42"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'7.220"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'159.72"
$ws.Range("D25").Value = "'8.521"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'0.1387"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'1.491"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'4.106"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'4.019"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'0.05253"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'1.844"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'0.7373"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").Value = "1.243.61"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").Value = "'2.761"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "'0.01786"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'6.369"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'0.8998"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'102.01"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "1.983.01"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").Value = "'64.33"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").Value = "'0.5109"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "'0.4043"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").Value = "'8.886"
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'6.706"
$ws.Range("E51").Value = "  +0.22%  "
